$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.187.47"
$ws.Range("E2").Value = "  +1.84%  "

$ws.Range("D3").Value = "2.383.82"
$ws.Range("E3").Value = "  +4.16%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "303.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.03%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.02%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.510"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.90%  "

$ws.Range("E8").Value = "  -0.10%  "

$ws.Range("E9").Value = "  +2.62%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.26"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.10%  "

$ws.Range("E11").Value = "  +1.54%  "

$ws.Range("E12").Value = "  +2.36%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.44"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.08%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.81"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.03%  "

$ws.Range("D15").Value = "2.753.11"
$ws.Range("E15").Value = "  +3.84%  "

$ws.Range("D16").Value = "2.387.67"
$ws.Range("E16").Value = "  +3.51%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.810"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.49%  "

$ws.Range("D18").Value = "43.171.14"
$ws.Range("E18").Value = "  +1.90%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.20"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.73%  "

$ws.Range("E20").Value = "  +6.06%  "

$ws.Range("D21").Value = "0.0₃0890"
$ws.Range("E21").Value = "  +0.65%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.52"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.74%  "

$ws.Range("E23").Value = "  +1.20%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "235.39"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.11%  "

$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.02%  "

$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.45"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.53%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.92"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.40%  "

$ws.Range("E28").Value = "  +0.29%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.15"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.86%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.67"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.54%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.08%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.10"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.23%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0740"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.26%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.23"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.58%  "

$ws.Range("E35").Value = "  +7.47%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.104"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.31%  "

$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.30"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.60%  "

$ws.Range("B38").Value = "WEMIXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.29"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.40%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.80"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.18%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "22.49"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +13.14%  "

$ws.Range("E41").Value = "  +0.84%  "

$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "1.958.99"
$ws.Range("E42").Value = "  +0.73%  "

$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "102.94"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -37.23%  "

$ws.Range("E44").Value = "  +1.41%  "

$ws.Range("E45").Value = "  +2.35%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.76"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.72%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.21"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -10.55%  "

$ws.Range("D48").Value = "2.605.72"
$ws.Range("E48").Value = "  +3.51%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "52.83"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.66%  "

$ws.Range("E50").Value = "  +3.46%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.96"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.41%  "
